$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.786.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.227.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.60%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.224.64'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.77'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.389'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.791.88'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.43%  '

$ws.Range("E14").Value = '  -3.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.889.58'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.235.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000159'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '416.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.206'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.495'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.06%  '

$ws.Range("E28").Value = '  -1.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.35%  '

$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.84%  '

$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.99'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.40%  '

$ws.Range("E36").Value = '  -1.99%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.834.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.726'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.28%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0630'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.21%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.18%  '

$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '300.99'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0263'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.101'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.96%  '
